$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the two footer (signature) rows down by one: 25->26, 26->27 ---
# (copy formats only - bottom-up order so sources aren't clobbered before being read)
$ws.Range("B26:C26").Copy()
$ws.Range("B27:C27").PasteSpecial(-4122)
$ws.Range("H26:J26").Copy()
$ws.Range("H27:J27").PasteSpecial(-4122)

$ws.Range("B25:C25").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("H25:J25").Copy()
$ws.Range("H26:J26").PasteSpecial(-4122)

# old row 25 no longer exists as a separate row once its content has moved to row 26
$ws.Range("B25:C25").UnMerge()
$ws.Range("H25:J25").UnMerge()
$ws.Range("B25:J25").Clear()
$ws.Range("B27:C27").Merge()
$ws.Range("H27:J27").Merge()

# --- Data rows: add one more "Periodo Mora" row (6 periods instead of 5) ---
# Old last row (20, special bottom-border style) moves to row 21
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# New row 20 takes on the "normal" data-row style (same as rows 16-19)
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Footer text (content unchanged, just re-asserted after the format-only paste) ---
$ws.Range("B26").Value = "___________________________________"
$ws.Range("H26").Value = "___________________________________"
$ws.Range("B27").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H27").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# --- Data rows 16-21: Tipo Doc / No Doc / Nombre / Periodo Mora / Valor Mora / Salario Basico ---
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047364257"
$ws.Range("D16").Value = "DARLYN ELENA NUÑEZ FUENTES"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047364257"
$ws.Range("D17").Value = "DARLYN ELENA NUÑEZ FUENTES"
$ws.Range("E17").Value = "2503"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047364257"
$ws.Range("D18").Value = "DARLYN ELENA NUÑEZ FUENTES"
$ws.Range("E18").Value = "2504"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047364257"
$ws.Range("D19").Value = "DARLYN ELENA NUÑEZ FUENTES"
$ws.Range("E19").Value = "2505"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047364257"
$ws.Range("D20").Value = "DARLYN ELENA NUÑEZ FUENTES"
$ws.Range("E20").Value = "2506"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1047364257"
$ws.Range("D21").Value = "DARLYN ELENA NUÑEZ FUENTES"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# --- Valor Mora total grows by one more period (284700 + 56940) ---
$ws.Range("E11").Value = 341640
